# Salaries and Tasks Update
# Fills in the salary/task tracking sheet with this team's actual data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info -----------------------------------------------------
$ws.Range("B3").Value = "11/12/2020"          # Date

# --- Team member names + salary distribution --------------------------
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("B8").Value = 110

$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("B9").Value = 90

$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 90

$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 110

# Only 4 team members this week - clear the 5th member slot
$ws.Range("A12").ClearContents()

$ws.Range("B4").Value = "Limette"             # Team Name
$ws.Range("B5").Value = 4                     # Total Number of Team Members

# --- Shrink the task-header row now that the instructions are shorter -
$ws.Rows.Item(18).RowHeight = 39

# --- Tasks ------------------------------------------------------------
$ws.Range("A19").Value = "High Fidelity Prototype"
$ws.Range("B19").Value = "High Fidelity Prototype"

# --- Restore the selection Excel last had on save ----------------------
$ws.Range("F18").Select()
